# Apply updated cryptocurrency price/volume data to Sheet1.
# Values that look numeric (e.g. "1.004", "0.000008445") must be stored as
# literal text to match the source data (inline strings), not auto-converted
# to numbers by Excel's smart input parsing (which would also mangle
# trailing zeros and flip tiny values into scientific notation). We force
# text storage by switching NumberFormat to "@" right before the write, then
# restoring the cell's style back to "Normal" so no formatting changes persist.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.686.64"
$ws.Range("E2").Value = "  -6.77%  "
$ws.Range("D3").Value = "1.698.80"
$ws.Range("E3").Value = "  -5.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5145"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -12.69%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2649"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "22.23"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06302"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07358"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.05%  "
$ws.Range("D12").Value = "1.698.76"
$ws.Range("E12").Value = "  -5.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.522"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5811"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.00%  "
$ws.Range("D15").Value = "1.929.81"
$ws.Range("E15").Value = "  -5.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008445"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -12.89%  "
$ws.Range("D18").Value = "26.698.31"
$ws.Range("E18").Value = "  -6.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.017"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.005"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("E21").Value = "  -4.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "187.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -10.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.267"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.522"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1163"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.354"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05667"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.335"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.511"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.490"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.642"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.023"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6030"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.357"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.690"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.102.91"
$ws.Range("E39").Value = "  -3.93%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01614"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8587"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.853"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.003"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").Value = "1.857.69"
$ws.Range("E45").Value = "  -4.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000111"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.157"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05245"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4320"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.50%  "
